# 自动更新Excel文件 - 2025-11-04 23:14:43
# Recompute the "剩余" (days remaining) column (E) from "总天" (total days, D)
# and "开始时间" (start date, F) as of "today". When a row's countdown has
# run out (remaining <= 0), the cycle restarts: remaining resets to the
# total and the start date is bumped to "today".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference "today" used for this run's recalculation.
$todayY = 2025
$todayM = 11
$todayD = 5

function Get-DayNumber($y, $m, $d) {
    # Julian-day-number style serial so we can diff two yyyymmdd dates with
    # plain integer math (no reliance on DateTime/TimeSpan COM plumbing).
    if ($m -le 2) {
        $y = $y - 1
        $m = $m + 12
    }
    $a = [Math]::Floor($y / 100)
    $b = 2 - $a + [Math]::Floor($a / 4)
    return [Math]::Floor(365.25 * ($y + 4716)) + [Math]::Floor(30.6001 * ($m + 1)) + $d + $b - 1524
}

$todaySerial = Get-DayNumber $todayY $todayM $todayD
$todayYmd = $todayY * 10000 + $todayM * 100 + $todayD

# Find the last used row in the sheet (data starts at row 2, header at row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2
    $startRaw = $ws.Cells.Item($r, 6).Value2

    if ($total -eq $null -or $startRaw -eq $null) { continue }

    $startStr = [string][int64]$startRaw
    if ($startStr.Length -ne 8) {
        # Malformed start date (e.g. "202510929") - leave the row untouched.
        continue
    }

    $sy = [int]$startStr.Substring(0, 4)
    $sm = [int]$startStr.Substring(4, 2)
    $sd = [int]$startStr.Substring(6, 2)

    $startSerial = Get-DayNumber $sy $sm $sd
    $elapsed = $todaySerial - $startSerial
    $remaining = $total - $elapsed

    if ($remaining -le 0) {
        # Cycle finished - restart it as of today.
        $ws.Cells.Item($r, 5).Value = $total
        $ws.Cells.Item($r, 6).Value = $todayYmd
    } else {
        $ws.Cells.Item($r, 5).Value = $remaining
    }
}
